$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The score columns (home_score / away_score) used a slightly different
# font (theme-colored) than the rest of the sheet; normalize them to the
# same explicit black font used elsewhere in the table.
$ws.Range("E2:F52").Font.Color = 0

# Fill in results for the matches that have since been played.
$ws.Range("E9").Value = 0   # Slovenia 0
$ws.Range("F9").Value = 1   # Denmark 1

$ws.Range("E10").Value = 0  # Serbia 0
$ws.Range("F10").Value = 1  # England 1

$ws.Range("E11").Value = 3  # Poland 3
$ws.Range("F11").Value = 0  # Netherlands 0
